# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 43   # 合肥·之心城购物中心-2024漫趣地带嘉年华（免费）: 42 -> 43
$wsExpo.Range("F6").Value = 5308 # 合肥·第九届环形宇宙动漫游戏嘉年华: 5303 -> 5308
$wsExpo.Range("F8").Value = 109  # 合肥·九号幻想动漫游戏嘉年华: 104 -> 109
$wsExpo.Range("F10").Value = 364 # 合肥·心动恋章·冬日序国乙&代号鸢同人only: 363 -> 364

# Sheet "全部类型" (All types) mirrors the same events at different rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 43    # 合肥·之心城购物中心-2024漫趣地带嘉年华（免费）: 42 -> 43
$wsAll.Range("F9").Value = 5308  # 合肥·第九届环形宇宙动漫游戏嘉年华: 5303 -> 5308
$wsAll.Range("F11").Value = 109  # 合肥·九号幻想动漫游戏嘉年华: 104 -> 109
$wsAll.Range("F14").Value = 364  # 合肥·心动恋章·冬日序国乙&代号鸢同人only: 363 -> 364
